$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2304.3
$ws.Range("J112").Value = 2999.8
$ws.Range("L112").Value = 8999.400000000001
$ws.Range("N112").Value = -11215.4

$ws.Range("H125").Value = 320
$ws.Range("I125").Value = 320
$ws.Range("K125").Value = 2880
$ws.Range("M125").Value = -420

$ws.Range("H132").Value = 2586.7693
$ws.Range("I132").Value = 1136.5555
$ws.Range("J132").Value = 5849.75
$ws.Range("K132").Value = 3409.6665
$ws.Range("L132").Value = 17549.25
$ws.Range("M132").Value = -879.6664999999998
$ws.Range("N132").Value = -22609.25

$ws.Range("H137").Value = 2919.4736
$ws.Range("I137").Value = 1779.125
$ws.Range("J137").Value = 3748.818
$ws.Range("K137").Value = 5337.375
$ws.Range("L137").Value = 11246.454
$ws.Range("M137").Value = -2787.375
$ws.Range("N137").Value = -16346.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 2997.5
$ws.Range("J10").Value = 2997.5
$ws.Range("L10").Value = 2997.5
$ws.Range("N10").Value = -3337.5

$ws.Range("H23").Value = 25007
$ws.Range("J23").Value = 25007
$ws.Range("L23").Value = 25007
$ws.Range("N23").Value = -25525

$ws.Range("H32").Value = 6989.8335
$ws.Range("I32").Value = 6989.8335
$ws.Range("K32").Value = 6989.8335
$ws.Range("M32").Value = -6702.8335

$ws.Range("H132").Value = 3941
$ws.Range("I132").Value = 3853
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 11559
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -9029
$ws.Range("N132").Value = -17059.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 416.33334
$ws.Range("I7").Value = 562
$ws.Range("K7").Value = 562
$ws.Range("M7").Value = -449

$ws.Range("H33").Value = 49999.668
$ws.Range("I33").Value = 39999.5
$ws.Range("J33").Value = 54999.75
$ws.Range("K33").Value = 39999.5
$ws.Range("L33").Value = 54999.75
$ws.Range("M33").Value = -39663.5
$ws.Range("N33").Value = -55671.75

$ws.Range("H86").Value = 5649.4
$ws.Range("I86").Value = 3833
$ws.Range("J86").Value = 8374
$ws.Range("K86").Value = 3833
$ws.Range("L86").Value = 8374
$ws.Range("M86").Value = -2710
$ws.Range("N86").Value = -10620

$ws.Range("H89").Value = 5649.4
$ws.Range("I89").Value = 3833
$ws.Range("J89").Value = 8374
$ws.Range("K89").Value = 19165
$ws.Range("L89").Value = 41870
$ws.Range("M89").Value = -13549
$ws.Range("N89").Value = -53102

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 40000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 40000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 40000
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -40340

$ws.Range("H35").Value = 12532.6
$ws.Range("I35").Value = 3165.75
$ws.Range("J35").Value = 50000
$ws.Range("K35").Value = 3165.75
$ws.Range("L35").Value = 50000
$ws.Range("M35").Value = -2871.75
$ws.Range("N35").Value = -50588

$ws.Range("H38").Value = 34600
$ws.Range("I38").Value = 3800
$ws.Range("J38").Value = 50000
$ws.Range("K38").Value = 3800
$ws.Range("L38").Value = 50000
$ws.Range("M38").Value = -3423
$ws.Range("N38").Value = -50754

$ws.Range("H46").Value = 34600
$ws.Range("I46").Value = 3800
$ws.Range("J46").Value = 50000
$ws.Range("K46").Value = 3800
$ws.Range("L46").Value = 50000
$ws.Range("M46").Value = -3589
$ws.Range("N46").Value = -50422

$ws.Range("H58").Value = 3635.8
$ws.Range("I58").Value = 3821.625
$ws.Range("J58").Value = 2892.5
$ws.Range("K58").Value = 3821.625
$ws.Range("L58").Value = 2892.5
$ws.Range("M58").Value = -3618.625
$ws.Range("N58").Value = -3298.5

$ws.Range("H99").Value = 2837.1428
$ws.Range("I99").Value = 2837.1428
$ws.Range("K99").Value = 2837.1428
$ws.Range("M99").Value = -1339.1428

$ws.Range("H126").Value = 2837.1428
$ws.Range("I126").Value = 2837.1428
$ws.Range("K126").Value = 8511.428400000001
$ws.Range("M126").Value = -6041.428400000001

$ws.Range("H132").Value = 115033.11
$ws.Range("J132").Value = 6374.5
$ws.Range("L132").Value = 19123.5
$ws.Range("N132").Value = -24183.5

$ws.Range("H136").Value = 3635.8
$ws.Range("I136").Value = 3821.625
$ws.Range("J136").Value = 2892.5
$ws.Range("K136").Value = 11464.875
$ws.Range("L136").Value = 8677.5
$ws.Range("M136").Value = -8914.875
$ws.Range("N136").Value = -13777.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3191.25
$ws.Range("I4").Value = 1298.3334
$ws.Range("J4").Value = 3822.2222
$ws.Range("K4").Value = 3895.0002
$ws.Range("L4").Value = 11466.6666
$ws.Range("M4").Value = -3783.0002
$ws.Range("N4").Value = -11690.6666

$ws.Range("H5").Value = 1195.409
$ws.Range("I5").Value = 1079.9
$ws.Range("K5").Value = 3239.7
$ws.Range("M5").Value = -3127.7

$ws.Range("H14").Value = 314.9
$ws.Range("I14").Value = 314.9
$ws.Range("K14").Value = 944.6999999999999
$ws.Range("M14").Value = -771.6999999999999

$ws.Range("H92").Value = 750

$ws.Range("H95").Value = 7500
$ws.Range("J95").Value = 7500
$ws.Range("L95").Value = 22500
$ws.Range("N95").Value = -26618

$ws.Range("H131").Value = 1645.8334
$ws.Range("J131").Value = 1687.5
$ws.Range("L131").Value = 5062.5
$ws.Range("N131").Value = -15142.5

$ws.Range("H135").Value = 1195.409
$ws.Range("I135").Value = 1079.9
$ws.Range("K135").Value = 9719.1
$ws.Range("M135").Value = -7184.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2054.077
$ws.Range("I113").Value = 2044
$ws.Range("J113").Value = 2065.8333
$ws.Range("K113").Value = 2044
$ws.Range("L113").Value = 2065.8333
$ws.Range("M113").Value = 126
$ws.Range("N113").Value = -6405.8333

$ws.Range("H126").Value = 3991.5
$ws.Range("I126").Value = 3985
$ws.Range("K126").Value = 11955
$ws.Range("M126").Value = -9485

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2707.9092
$ws.Range("I82").Value = 1650.5
$ws.Range("J82").Value = 5527.6665
$ws.Range("K82").Value = 1650.5
$ws.Range("L82").Value = 5527.6665
$ws.Range("M82").Value = -1289.5
$ws.Range("N82").Value = -6249.6665

$ws.Range("H85").Value = 2707.9092
$ws.Range("I85").Value = 1650.5
$ws.Range("J85").Value = 5527.6665
$ws.Range("K85").Value = 1650.5
$ws.Range("L85").Value = 5527.6665
$ws.Range("M85").Value = -402.5
$ws.Range("N85").Value = -8023.6665

$ws.Range("H93").Value = 1551
$ws.Range("I93").Value = 1551
$ws.Range("K93").Value = 1551
$ws.Range("M93").Value = -303

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 17837832
$ws.Range("I2").Value = 20405400
$ws.Range("K2").Value = 20405400
$ws.Range("M2").Value = -20405288

$ws.Range("H4").Value = 1318724
$ws.Range("J4").Value = 1684933.4
$ws.Range("L4").Value = 1684933.4
$ws.Range("N4").Value = -1685159.4

$ws.Range("H40").Value = 3025
$ws.Range("I40").Value = 3025
$ws.Range("K40").Value = 3025
$ws.Range("M40").Value = -2876
